$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "#muliassus"
$ws.Range("C2").Value = "Muliassus"
$ws.Range("B3").Value = "#eunuchus.-bordello"
$ws.Range("C3").Value = "Eunuchus. Bordello"
$ws.Range("B4").Value = "#mu"
$ws.Range("C4").Value = "Mu"
$ws.Range("B5").Value = "#themoclea"
$ws.Range("C5").Value = "Themoclea"
$ws.Range("B6").Value = "#raad"
$ws.Range("C6").Value = "Raad"
$ws.Range("B7").Value = "#eunuchus"
$ws.Range("C7").Value = "Eunuchus"
$ws.Range("B8").Value = "#ferrando,-half-gekleedt,-wandelende-door-zijn-kamer"
$ws.Range("C8").Value = "Ferrando, half gekleedt, wandelende door zijn kamer"
$ws.Range("B9").Value = "#iulia-en-amada"
$ws.Range("C9").Value = "Iulia en Amada"
$ws.Range("B10").Value = "#bor"
$ws.Range("C10").Value = "Bor"
$ws.Range("B11").Value = "#feni"
$ws.Range("C11").Value = "Feni"
$ws.Range("B12").Value = "#them"
$ws.Range("C12").Value = "Them"
$ws.Range("B13").Value = "#bordel"
$ws.Range("C13").Value = "Bordel"
$ws.Range("B14").Value = "#borgias"
$ws.Range("C14").Value = "Borgias"
$ws.Range("B15").Value = "#fen"
$ws.Range("C15").Value = "Fen"
$ws.Range("B16").Value = "#per"
$ws.Range("C16").Value = "Per"
$ws.Range("B17").Value = "#the"
$ws.Range("C17").Value = "The"
$ws.Range("B18").Value = "#themoclea,-als-een-geest"
$ws.Range("C18").Value = "Themoclea, als een Geest"
$ws.Range("B19").Value = "#ferran"
$ws.Range("C19").Value = "Ferran"
$ws.Range("B20").Value = "#phego"
$ws.Range("C20").Value = "Phego"
$ws.Range("B21").Value = "#mul"
$ws.Range("C21").Value = "Mul"
$ws.Range("B22").Value = "#phil"
$ws.Range("C22").Value = "Phil"
$ws.Range("B23").Value = "#borgias.-themoclea.-muliassus"
$ws.Range("C23").Value = "Borgias. Themoclea. Muliassus"
$ws.Range("B24").Value = "#ama"
$ws.Range("C24").Value = "Ama"
$ws.Range("B25").Value = "#fer"
$ws.Range("C25").Value = "Fer"
$ws.Range("B26").Value = "#amad"
$ws.Range("C26").Value = "Amad"
$ws.Range("B27").Value = "#borgias-antwoordt-van-binnen"
$ws.Range("C27").Value = "Borgias antwoordt van binnen"
$ws.Range("B28").Value = "#iulia.-muliassus"
$ws.Range("C28").Value = "Iulia. Muliassus"
$ws.Range("B29").Value = "#iul"
$ws.Range("C29").Value = "Iul"
$ws.Range("B30").Value = "#borg"
$ws.Range("C30").Value = "Borg"
$ws.Range("B31").Value = "#eun"
$ws.Range("C31").Value = "Eun"
$ws.Range("B32").Value = "#fenizo-op-'t-ledekant.-themoclea"
$ws.Range("C32").Value = "Fenizo op 't Ledekant. Themoclea"

$ws.Range("D2:D23").ClearContents()
